# Repull data, push all data, mean calculation
# Updates column F (dSF) values on the active worksheet to reflect
# the repulled/recalculated data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -3
    3  = 1
    4  = 1
    5  = -1
    7  = -2
    8  = 9
    9  = 2
    10 = -1
    11 = 5
    12 = -3
    14 = 1
    17 = 7
    19 = -3
    20 = -1
    21 = -1
    22 = 1
    23 = 1
    24 = 1
    25 = 4
    26 = 1
    27 = -1
    28 = -2
    29 = 3
    30 = 4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
